$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheets 1-4: "Potencia Acumulada - SIN (MW)", "Geracao Periodo Medio (MWMed)",
#             "Atendimento a Ponta(MW)", "Potencia Incremental - SIN(MW)"
# Same transformation applied to each:
#   - Add header cell A1 = "Fonte/Tecnologia" (styled like B1:E1)
#   - Remove the bold/border style from A2:A12 (row-label cells)
#   - Fix a few accented labels
# ---------------------------------------------------------------------------
foreach ($i in 1..4) {
    $ws = $wb.Worksheets.Item($i)

    # Add the new header cell, copying the header-row formatting from B1.
    $ws.Range("A1").Value = "Fonte/Tecnologia"
    $ws.Range("B1").Copy()
    $ws.Range("A1").PasteSpecial(-4122)
    $excel.CutCopyMode = $false

    # Update row-label text where accents were added/fixed.
    $ws.Range("A3").Value = "Gás Natural"
    $ws.Range("A4").Value = "Carvão"
    $ws.Range("A6").Value = "Óleos Comb"
    $ws.Range("A8").Value = "Eólica"
    $ws.Range("A11").Value = "Pot. Compl."

    # Remove the bold/border styling from all the row-label cells A2:A12.
    $ws.Range("A2:A12").ClearFormats()
}

# ---------------------------------------------------------------------------
# Sheet 5: "Emissoes Totais (MtCO2eq)"
#   - Add header cell A1 = "Período" (styled like B1:E1)
#   - Rename A2/A3 labels, drop their bold style
#   - Remove row 4 ("Teto") entirely
# ---------------------------------------------------------------------------
$ws5 = $wb.Worksheets.Item(5)

$ws5.Range("A1").Value = "Período"
$ws5.Range("B1").Copy()
$ws5.Range("A1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws5.Range("A2").Value = "P.Médio"
$ws5.Range("A3").Value = "P.Crítico"
$ws5.Range("A2:A3").ClearFormats()

$ws5.Rows("4:4").Delete()

# ---------------------------------------------------------------------------
# Sheet 6: "Custo Total (bilhões de R$)"
#   - Add header cell A1 = "Tipo Expansão" (styled like B1)
#   - Change B1 label from "Custo" to "2015"
#   - Rename A2/A3 labels, drop their bold style
#   - Update B2/B3 values
# ---------------------------------------------------------------------------
$ws6 = $wb.Worksheets.Item(6)

# B1 must stay a *text* "2015" (like the other sheets' header cells), not a
# number. A plain Value assignment would auto-coerce the digit string to a
# number, so build it as a text formula in a scratch cell (guaranteeing a
# string result) and paste only the value back into B1, preserving B1's
# existing (bold/border) style untouched.
$tmp6 = $ws6.Range("Z100")
$tmp6.Formula = "=""20""&""15"""
$tmp6.Copy()
$ws6.Range("B1").PasteSpecial(-4163)
$excel.CutCopyMode = $false
$tmp6.Delete()

$ws6.Range("A1").Value = "Tipo Expansão"
$ws6.Range("B1").Copy()
$ws6.Range("A1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws6.Range("A2").Value = "Expansão Centralizada"
$ws6.Range("B2").Value = 571
$ws6.Range("A3").Value = "Expansão por GD"
$ws6.Range("B3").Value = 99
$ws6.Range("A2:A3").ClearFormats()
